$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as Text in the source sheet even when their
# content looks numeric (e.g. "305.12", "1.633.84"). Assigning such a literal to
# Range.Value directly would make Excel auto-convert it to a Number, so we use a
# leading apostrophe to force text entry, then reset Style to Normal so we don't
# leave a stray quote-prefixed style on the cell (matches original formatting).
$ws.Range("D2").Value = "'23.425.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "'1.639.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'305.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").Value = "'0.3728"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("D8").Value = "'52.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.3618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "'1.252"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").Value = "'0.08120"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'22.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'6.590"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "'0.00001268"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "'7.283"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "'1.633.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "'94.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'0.06879"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'18.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'6.512"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'23.425.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").Value = "'12.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'3.036"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.405"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "'21.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'151.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "'5.311"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").Value = "'135.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").Value = "'2.293"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("D32").Value = "'1.813.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").Value = "'6.745"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "'0.9530"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "'0.02843"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.07%  "
$ws.Range("D36").Value = "'10.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").Value = "'0.2515"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'0.07232"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("D39").Value = "'0.08775"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'6.049"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "'1.372"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'0.7043"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").Value = "'12.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "'16.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "'0.6507"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.325"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "'4.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").Value = "'0.07967"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").Value = "'128.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").Value = "'1.200"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.44%  "
